# Update the public EPEX Spot prices workbook:
# add a new column AT ("29-jul") to the "Prix Spot" sheet, mirroring the
# existing header style from column AS and filling in the daily prices.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Prix Spot")

# New header cell: copy the formatting (border/bold/alignment) of the
# previous header cell (AS1) so the new column matches the existing
# header style, then set its text.
$ws.Range("AS1").Copy()
$ws.Range("AT1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("AT1").Value = "29-jul"

# New data values for the "29-jul" column (rows 2-25).
$ws.Range("AT2").Value = 50.68
$ws.Range("AT3").Value = 37.97
$ws.Range("AT4").Value = 34.78
$ws.Range("AT5").Value = 32
$ws.Range("AT6").Value = 22.98
$ws.Range("AT7").Value = 34.78
$ws.Range("AT8").Value = 45
$ws.Range("AT9").Value = 50
$ws.Range("AT10").Value = 49.45
$ws.Range("AT11").Value = 35.61
$ws.Range("AT12").Value = 11.73
$ws.Range("AT13").Value = 25.25
$ws.Range("AT14").Value = 25.92
$ws.Range("AT15").Value = 5.34
$ws.Range("AT16").Value = 3.78
$ws.Range("AT17").Value = 8.51
$ws.Range("AT18").Value = 22.87
$ws.Range("AT19").Value = 51.53
$ws.Range("AT20").Value = 48.97
$ws.Range("AT21").Value = 72.14
$ws.Range("AT22").Value = 57.58
$ws.Range("AT23").Value = 98
$ws.Range("AT24").Value = 103.27
$ws.Range("AT25").Value = 84.13

Write-Output "AT column (29-jul) added to Prix Spot sheet"
